$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row (B1:E1)
$ws.Cells.Item(1,2).Value = "Form"
$ws.Cells.Item(1,3).Value = "Goals scored"
$ws.Cells.Item(1,4).Value = "Goals conceded"
$ws.Cells.Item(1,5).Value = "Total Goals"

# Column A rows 2-24: row index "1".."23" as TEXT (matches existing shared
# strings "1".."23" used elsewhere in the workbook). Writing a ="n" formula
# then pasting values-only keeps the text type without adding a new style.
for ($i = 2; $i -le 24; $i++) {
    $n = $i - 1
    $ws.Cells.Item($i, 1).Formula = "=""$n"""
}
$ws.Range("A2:A24").Copy()
$ws.Range("A2:A24").PasteSpecial(-4163)


# Column B: Form strings, one per team, in team order (rows 2-24)
$form = @(
    "Aldershot,W L D W L L"
    "Altrincham,L L L L D L"
    "Barnet,W D W L L W"
    "Boreham Wood,W D L D L W"
    "Bromley,W W W W L W"
    "Chesterfield,L L D W L L"
    "Dag and Red,L W W W W W"
    "Dover Athletic,L D L L W L"
    "Eastleigh,D L W W W L"
    "Halifax,W W L W L W"
    "Hartlepool,W D W W W L"
    "Kings Lynn,L D L L D L"
    "Maidenhead,D W W L D D"
    "Notts County,L L L D L D"
    "Solihull,W L W L L W"
    "Stockport,W D W W W W"
    "Sutton,D L W L W W"
    "Torquay,W W D W W W"
    "Wealdstone,L L L L L L"
    "Weymouth,L L L W W L"
    "Woking,L L L L L L"
    "Wrexham,L W W D W D"
    "Yeovil,W W L W W L"
)
for ($i = 0; $i -lt $form.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $form[$i]
}

# Column C: Goals scored strings, one per team, in team order (rows 2-24)
$goalsScored = @(
    "Aldershot,3 1 1 2 1 1"
    "Altrincham,0 0 0 2 1 0"
    "Barnet,3 0 3 0 0 2"
    "Boreham Wood,3 0 0 0 0 5"
    "Bromley,1 4 2 2 1 1"
    "Chesterfield,1 1 0 2 1 0"
    "Dag and Red,0 2 3 3 2 3"
    "Dover Athletic,0 1 1 0 3 1"
    "Eastleigh,0 0 1 1 2 1"
    "Halifax,4 2 1 1 0 3"
    "Hartlepool,2 2 2 7 3 0"
    "Kings Lynn,2 1 0 0 0 1"
    "Maidenhead,2 6 2 0 0 2"
    "Notts County,0 1 0 2 0 1"
    "Solihull,2 2 5 0 1 4"
    "Stockport,4 2 3 4 1 4"
    "Sutton,0 0 4 0 1 3"
    "Torquay,1 2 2 2 3 2"
    "Wealdstone,1 2 0 0 0 0"
    "Weymouth,2 1 0 1 4 0"
    "Woking,0 0 1 0 2 1"
    "Wrexham,0 4 4 0 2 2"
    "Yeovil,3 1 1 3 2 0"
)
for ($i = 0; $i -lt $goalsScored.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $goalsScored[$i]
}

# Column D: Goals conceded strings, one per team, in team order (rows 2-24)
$goalsConceded = @(
    "Aldershot,0 3 1 0 2 3"
    "Altrincham,2 1 4 3 1 4"
    "Barnet,1 0 2 1 2 0"
    "Boreham Wood,1 0 1 0 3 1"
    "Bromley,0 3 1 1 2 0"
    "Chesterfield,2 2 0 1 3 2"
    "Dag and Red,1 0 2 0 0 1"
    "Dover Athletic,1 1 3 2 1 3"
    "Eastleigh,0 2 0 0 0 3"
    "Halifax,2 1 2 0 1 0"
    "Hartlepool,1 2 0 2 1 1"
    "Kings Lynn,4 1 3 4 0 5"
    "Maidenhead,2 0 1 1 0 2"
    "Notts County,2 2 1 2 2 1"
    "Solihull,1 3 1 3 2 0"
    "Stockport,0 2 0 0 0 0"
    "Sutton,0 1 0 1 0 1"
    "Torquay,0 0 2 1 1 0"
    "Wealdstone,3 7 6 2 2 4"
    "Weymouth,3 2 3 0 2 2"
    "Woking,4 2 2 1 4 3"
    "Wrexham,3 0 0 0 1 2"
    "Yeovil,1 0 5 0 0 3"
)
for ($i = 0; $i -lt $goalsConceded.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $goalsConceded[$i]
}

# Column E: Total Goals strings, one per team, in team order (rows 2-24)
$totalGoals = @(
    "Aldershot,3 4 2 2 3 4"
    "Altrincham,2 1 4 5 2 4"
    "Barnet,4 0 5 1 2 2"
    "Boreham Wood,4 0 1 0 3 6"
    "Bromley,1 7 3 3 3 1"
    "Chesterfield,3 3 0 3 4 2"
    "Dag and Red,1 2 5 3 2 4"
    "Dover Athletic,1 2 4 2 4 4"
    "Eastleigh,0 2 1 1 2 4"
    "Halifax,6 3 3 1 1 3"
    "Hartlepool,3 4 2 9 4 1"
    "Kings Lynn,6 2 3 4 0 6"
    "Maidenhead,4 6 3 1 0 4"
    "Notts County,2 3 1 4 2 2"
    "Solihull,3 5 6 3 3 4"
    "Stockport,4 4 3 4 1 4"
    "Sutton,0 1 4 1 1 4"
    "Torquay,1 2 4 3 4 2"
    "Wealdstone,4 9 6 2 2 4"
    "Weymouth,5 3 3 1 6 2"
    "Woking,4 2 3 1 6 4"
    "Wrexham,3 4 4 0 3 4"
    "Yeovil,4 1 6 3 2 3"
)
for ($i = 0; $i -lt $totalGoals.Count; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $totalGoals[$i]
}

$wb.Worksheets.Item(1).Activate()
